$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 2 (same values as row 1: admin / manager)
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Add row 3 (new values: mangh / chitgup)
$ws.Range("A3").Value = "mangh"
$ws.Range("B3").Value = "chitgup"

# Update the selection to C9, matching the target sheet's sheetView selection
$ws.Range("C9").Select()
